$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Jim Hadid"
$ws.Range("C12").Value = "jherkhdbg121@gmail.com"
$ws.Range("D12").Value = 9226475786

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Hemma Hadid"
$ws.Range("C13").Value = "jhfrgeywuhdsj121@gmail.com"
$ws.Range("D13").Value = 9226461104
